$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = "{'perception,sensing,socket': 4, 'perception,planning': 2, 'planning,sensing': 1, 'sensing,socket': 2, 'data,perception': 1, 'perception,sensing': 5}"
$ws.Range("H2").Value = 154

$ws.Range("G3").Value = "{'planning,sensing': 1, 'perception,planning': 10, 'perception,planning,sensing': 1, 'perception,planning,sensing,socket': 2, 'perception,sensing,socket': 1, 'data,perception,planning': 1, 'perception,sensing': 5, 'sensing,socket': 2}"
$ws.Range("H3").Value = 223

$ws.Range("G4").Value = "{'perception,planning': 22, 'perception,system': 1, 'planning,sensing': 6, 'perception,planning,socket': 4, 'perception,planning,sensing,socket': 2, 'sensing,socket': 1, 'perception,sensing': 5, 'data,planning': 3, 'perception,planning,sensing': 2}"
$ws.Range("H4").Value = 389

$ws.Range("G5").Value = "{'perception,sensing': 4, 'planning,sensing': 6, 'perception,planning,sensing': 4, 'perception,planning,sensing,socket': 5, 'perception,planning': 18, 'data,perception,planning,sensing': 1, 'perception,planning,socket': 1, 'sensing,socket': 3}"
$ws.Range("H5").Value = 370

$ws.Range("G6").Value = "{'perception,planning,sensing,socket': 5, 'perception,planning': 26, 'perception,planning,sensing': 4, 'data,perception,planning,sensing': 1, 'perception,sensing': 5, 'planning,sensing': 5, 'perception,planning,socket': 1, 'sensing,socket': 3, 'data,planning': 1, 'data,socket': 1}"
$ws.Range("H6").Value = 465

$ws.Range("G7").Value = "{'perception,planning': 29, 'perception,planning,sensing,socket,util': 2, 'planning,sensing': 5, 'perception,planning,sensing,socket': 3, 'perception,planning,sensing': 5, 'perception,planning,socket': 1, 'data,perception,planning,sensing': 1, 'perception,sensing': 3, 'perception,util': 2, 'data,planning': 2, 'data,perception,planning': 1, 'data,perception': 2, 'sensing,socket': 4, 'data,socket': 1, 'planning,system': 1, 'perception,system': 1}"
$ws.Range("H7").Value = 581

$ws.Range("G8").Value = "{'perception,planning': 9, 'perception,planning,sensing': 4, 'perception,planning,sensing,socket': 1, 'planning,sensing,socket': 1, 'planning,sensing': 11, 'perception,sensing,util': 1, 'perception,sensing,socket,util': 1, 'data,perception,planning,sensing': 1, 'data,planning': 1, 'perception,sensing': 4, 'data,socket': 1, 'data,perception,planning,socket': 1, 'actuation,perception,planning': 7, 'actuation,perception': 1, 'sensing,socket': 2, 'actuation,planning': 2, 'planning,system': 1, 'perception,util': 1, 'perception,planning,system': 1, 'data,perception': 2}"
$ws.Range("H8").Value = 513

$ws.Range("G9").Value = "{'perception,planning': 44, 'actuation,perception,planning': 7, 'perception,planning,sensing': 14, 'perception,sensing': 16, 'actuation,planning': 2, 'data,perception,planning,sensing': 1, 'perception,planning,sensing,socket': 2, 'perception,planning,util': 1, 'perception,planning,socket,util': 1, 'data,perception,planning': 2, 'planning,system': 2, 'perception,util': 1, 'planning,sensing': 4, 'actuation,perception': 2, 'data,planning': 2, 'perception,planning,system': 1, 'perception,system': 1, 'planning,sensing,socket': 1, 'data,perception': 2}"
$ws.Range("H9").Value = 1130
